$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 10000
$ws.Range("B2").Value = 1333

$ws.Range("B2").Select()

$wb.Save()
